$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 89, shifting existing rows 89:197 down to 90:198
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the new data record
$ws.Range("A89").Value = 3
$ws.Range("B89").Value = "Femacal de La Calera"
$ws.Range("C89").Value = "Coquimbo"
$ws.Range("D89").Value = 44546
$ws.Range("E89").Value = 5
$ws.Range("F89").Value = 100112001
$ws.Range("G89").Value = "Berenjena"
$ws.Range("H89").Value = "Sin especificar"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 110
$ws.Range("K89").Value = 8000
$ws.Range("L89").Value = 8500
$ws.Range("M89").Value = 8227
$ws.Range("N89").Value = "$/caja 60 unidades"
$ws.Range("O89").Value = "Región de Arica y Parinacota"
$ws.Range("P89").Value = 137
$ws.Range("Q89").Value = 60
$ws.Range("R89").Value = "Hortaliza"
